$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week1")

# Set "Actual time length to complete" for row 6 (Discussion question 2) = 1:08 ([h]:mm -> 1 hour 8 min)
$ws.Range("C6").Value = 0.047222222222222221

# Set "Actual time length to complete" for row 12 (Hand-in assignment) = 0:45 ([h]:mm -> 45 min)
$ws.Range("C12").Value = 0.03125

# Update the active selection to C13 as recorded in the saved workbook view
$ws.Range("C13").Select()
